# Lunggo_Config.xlsx edit
# Commit: "Mystifly Production Web (Impotent), Minor fix for Mystifly and API"
#
# Inserts a new "mystifly.apiEndPoint" row (new row 53) ahead of the existing
# mystifly block (apiAccountNumber/apiUserName/apiPassword/apiTargetServer,
# which shift down from rows 53-56 to 54-57), fills in the new Production
# web-service URLs/credentials for mystifly in column F, and flips row 57
# (apiTargetServer) column F from "Test" to "Production".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row above the "mystifly" block (old row 53) ---------------
$ws.Rows(53).Insert()

# --- New row 53: mystifly / apiEndPoint -------------------------------------
$ws.Range("A53").Value = "*"
$ws.Range("B53").Value = "mystifly"
$ws.Range("C53").Value = "apiEndPoint"
$ws.Range("C53").Style = "Normal"
$ws.Range("A53:C53").Style = "Normal"
$ws.Range("A53").Value = "*"
$ws.Range("B53").Value = "mystifly"
$ws.Range("C53").Value = "apiEndPoint"
$ws.Range("D53").Formula = '="@@."&A53&"."&B53&"."&C53&"@@"'
$ws.Range("E53").Value = "http://apidemo.myfarebox.com/V2/OnePoint.svc"
$ws.Range("G53").Value = "http://apidemo.myfarebox.com/V2/OnePoint.svc"
$ws.Range("H53").Value = ""
$ws.Range("I53").Value = ""

# Match the formatting used throughout the sheet for this table section
$ws.Range("A53:D53").Style = "Normal"
$ws.Range("E53:I53").Style = "Normal"

$ws.Range("A53").Value = "*"
$ws.Range("B53").Value = "mystifly"
$ws.Range("C53").Value = "apiEndPoint"
$ws.Range("D53").Formula = '="@@."&A53&"."&B53&"."&C53&"@@"'
$ws.Range("E53").Value = "http://apidemo.myfarebox.com/V2/OnePoint.svc"
$ws.Range("G53").Value = "http://apidemo.myfarebox.com/V2/OnePoint.svc"

$ws.Range("A53").Copy()
$ws.Range("A2").PasteSpecial()
